$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $rowIndex, $newText) {
    $cell = $table.Cell($rowIndex, 1)
    $r = $cell.Range
    $r.End = $r.End - 1
    $r.Text = $newText
}

# Single-value rows near the top of the table
Set-CellText $t 1 "0M"
Set-CellText $t 2 "0M"
Set-CellText $t 3 "0M"
Set-CellText $t 4 "303"
Set-CellText $t 5 "0.00003"
Set-CellText $t 6 "0.00044"
Set-CellText $t 7 "0.00015"
Set-CellText $t 9 "0.00021"
Set-CellText $t 10 "0.00024"
Set-CellText $t 11 "0.00029"
Set-CellText $t 12 "0.05128"

# Rows near the end of the table that currently contain tab-separated
# multi-run values get collapsed down to a single value.
Set-CellText $t 44 "99.95"
Set-CellText $t 45 "0.05"
Set-CellText $t 46 "99"
